$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contract rows to append below the existing data (rows 22-24).
# Column A dates and several numeric-looking columns (Qty, NSN, P/N) must be
# stored as plain text, matching the existing rows further up the sheet, so
# force text formatting on the target range before writing any values.
$rng = $ws.Range("A22:J24")
$rng.NumberFormat = "@"

$ws.Range("A22").Value = "3/26/2019"
$ws.Range("B22").Value = "SPE7MC-19-V-5974"
$ws.Range("C22").Value = "71"
$ws.Range("D22").Value = "$14,743.86 "
$ws.Range("E22").Value = "5999016427529"
$ws.Range("F22").Value = "DELAY LINE"
$ws.Range("G22").Value = "Data Delay Devices"
$ws.Range("H22").Value = "583R874H08"
$ws.Range("I22").Value = "CP"
$ws.Range("J22").Value = "2019 SEP 03"

$ws.Range("A23").Value = "3/26/2019"
$ws.Range("B23").Value = "SPE4A4-19-V-5178"
$ws.Range("C23").Value = "5"
$ws.Range("D23").Value = "$1,568.35 "
$ws.Range("E23").Value = "6150015755067"
$ws.Range("F23").Value = "CABLE ASSEMBLY,POWER,ELECTRICAL"
$ws.Range("G23").Value = "DITMCO"
$ws.Range("H23").Value = "114-04030-0008"
$ws.Range("I23").Value = "CP"
$ws.Range("J23").Value = "2019 JUL 24"

$ws.Range("A24").Value = "3/27/2019"
$ws.Range("B24").Value = "SPE7M5-19-P-5961"
$ws.Range("C24").Value = "1000"
$ws.Range("D24").Value = "$43,380.00"
$ws.Range("E24").Value = "5935016786940"
$ws.Range("F24").Value = "BACKSHELL,ELECTRICAL CONNECTOR"
$ws.Range("G24").Value = "Glenair"
$ws.Range("H24").Value = "445FS065NF24343"
$ws.Range("I24").Value = "M41"
$ws.Range("J24").Value = "2019 SEP 03"

# Restore the default "Normal" style so the new rows don't keep an explicit
# text-number-format style index (matches the unstyled cells used by the
# rest of the sheet's text rows).
$rng.Style = "Normal"
